{"js": "// Replace two-digit-by-two-digit multiplication expressions in the table\n// with the new set of expressions, per the commit diff.\nconst replacements = [\n  [\"96\u00d752=\", \"28\u00d761=\"],\n  [\"93\u00d722=\", \"36\u00d767=\"],\n  [\"46\u00d782=\", \"69\u00d711=\"],\n  [\"21\u00d737=\", \"72\u00d778=\"],\n  [\"65\u00d714=\", \"61\u00d798=\"],\n  [\"65\u00d760=\", \"90\u00d716=\"],\n  [\"55\u00d747=\", \"18\u00d734=\"],\n  [\"91\u00d724=\", \"50\u00d746=\"],\n  [\"19\u00d718=\", \"96\u00d789=\"],\n  [\"51\u00d773=\", \"21\u00d797=\"],\n  [\"39\u00d770=\", \"41\u00d727=\"],\n  [\"32\u00d736=\", \"84\u00d776=\"],\n  [\"60\u00d752=\", \"55\u00d742=\"],\n  [\"68\u00d735=\", \"16\u00d733=\"],\n  [\"23\u00d765=\", \"90\u00d771=\"],\n  [\"37\u00d765=\", \"17\u00d757=\"],\n  [\"44\u00d758=\", \"90\u00d795=\"],\n  [\"87\u00d739=\", \"37\u00d773=\"],\n  [\"87\u00d790=\", \"19\u00d760=\"],\n  [\"84\u00d773=\", \"27\u00d737=\"],\n  [\"42\u00d716=\", \"98\u00d738=\"],\n  [\"86\u00d736=\", \"50\u00d715=\"],\n  [\"19\u00d766=\", \"41\u00d730=\"],\n  [\"81\u00d728=\", \"56\u00d760=\"],\n  [\"62\u00d761=\", \"55\u00d785=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace two-digit-by-two-digit multiplication expressions in the table\n# with the new set of expressions, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"96\u00d752=\", \"28\u00d761=\"),\n    @(\"93\u00d722=\", \"36\u00d767=\"),\n    @(\"46\u00d782=\", \"69\u00d711=\"),\n    @(\"21\u00d737=\", \"72\u00d778=\"),\n    @(\"65\u00d714=\", \"61\u00d798=\"),\n    @(\"65\u00d760=\", \"90\u00d716=\"),\n    @(\"55\u00d747=\", \"18\u00d734=\"),\n    @(\"91\u00d724=\", \"50\u00d746=\"),\n    @(\"19\u00d718=\", \"96\u00d789=\"),\n    @(\"51\u00d773=\", \"21\u00d797=\"),\n    @(\"39\u00d770=\", \"41\u00d727=\"),\n    @(\"32\u00d736=\", \"84\u00d776=\"),\n    @(\"60\u00d752=\", \"55\u00d742=\"),\n    @(\"68\u00d735=\", \"16\u00d733=\"),\n    @(\"23\u00d765=\", \"90\u00d771=\"),\n    @(\"37\u00d765=\", \"17\u00d757=\"),\n    @(\"44\u00d758=\", \"90\u00d795=\"),\n    @(\"87\u00d739=\", \"37\u00d773=\"),\n    @(\"87\u00d790=\", \"19\u00d760=\"),\n    @(\"84\u00d773=\", \"27\u00d737=\"),\n    @(\"42\u00d716=\", \"98\u00d738=\"),\n    @(\"86\u00d736=\", \"50\u00d715=\"),\n    @(\"19\u00d766=\", \"41\u00d730=\"),\n    @(\"81\u00d728=\", \"56\u00d760=\"),\n    @(\"62\u00d761=\", \"55\u00d785=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
